# Add the new "pd_method" assumption field (two Z-TM methods).
#
# 1. ASSUMPTIONS sheet: insert a new column D ("pd_method") populated with
#    the default "METHOD-1" for every segment row.
# 2. DICTIONARY sheet: insert a new row describing the "pd_method" field.
# 3. TRANSITION_MATRIX formulas that XLOOKUP into ASSUMPTIONS columns to
#    the right of the insertion point are automatically re-pointed by Excel
#    when the column is inserted (M->N, H->I, T->U, U->V).

$wb = $excel.ActiveWorkbook

$wsAssumptions = $wb.Worksheets.Item("ASSUMPTIONS")
$wsDictionary  = $wb.Worksheets.Item("DICTIONARY")

# --- 1. ASSUMPTIONS: insert new "pd_method" column before the current D ---
$oldColCWidth = $wsAssumptions.Columns("C:C").ColumnWidth
$wsAssumptions.Columns("D:D").Insert()
$wsAssumptions.Columns("D:D").ColumnWidth = $oldColCWidth

$wsAssumptions.Range("D1").Value = "pd_method"
$wsAssumptions.Range("D2:D5").Value = "METHOD-1"

# --- 2. DICTIONARY: insert a new row documenting "pd_method" ---
$wsDictionary.Rows("5:5").Insert()

# Copy the formatting of the row above (field row) onto the new row so the
# borders / alignment match the rest of the dictionary table.
$wsDictionary.Range("A4:C4").Copy()
$wsDictionary.Range("A5:C5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsDictionary.Range("A5").Value = "pd_method"
$wsDictionary.Range("B5").Value = "The Z-model method to use to calculate the FiT TM." + [char]10 + "METHOD-1: " + [char]10 + "    Use the Z-risk engine method where bin widths are also transformed via the Z-factor." + [char]10 + "    Typically leads to more macro sensitivity and higher PDs." + [char]10 + "METHOD-2:" + [char]10 + "    Use the Credit Metrics method where bin widths are fixed and only shifted by the Z-factor."
$wsDictionary.Range("C5").Value = "<string>"
$wsDictionary.Rows("5:5").RowHeight = 105

# --- 3. Restore/update the active selections to match the edited workbook ---
$wsAssumptions.Activate()
$wsAssumptions.Range("D8").Select()

$wsDictionary.Activate()
$wsDictionary.Range("B6").Select()

Write-Host "pd_method field added to ASSUMPTIONS and DICTIONARY"
